$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new test-case rows (17-25) to Tabelle1 ---
# Row 17
$ws.Range("A17").Value = "Insurant Data Page check for open mandatory fields"
$ws.Range("B17").Value = "<CHK>"
$ws.Range("D17").Value = "Check for open mandatory fields"
$ws.Range("H17").Value = "<NOP>"

# Row 18
$ws.Range("A18").Value = "103_TruckInsurance_003_InsurantData_001_MandatoryFields_FillFirstName"
$ws.Range("B18").Value = "<SET>"
$ws.Range("D18").Value = "MandatoryFields_FillFirstName"
$ws.Range("H18").Value = "<NOP>"

# Row 19
$ws.Range("A19").Value = "103_TruckInsurance_003_InsurantData_001_MandatoryFields_CheckFilledFirstName"
$ws.Range("B19").Value = "<CHK>"
$ws.Range("D19").Value = "MandatoryFields_CheckFilledFirstName"
$ws.Range("H19").Value = "<NOP>"

# Row 20
$ws.Range("A20").Value = "Insurant Data Page check for hints regarding mandatory fields"
$ws.Range("B20").Value = "<CHK>"
$ws.Range("D20").Value = "Check for hints regarding mandatory fields"
$ws.Range("H20").Value = "<NOP>"

# Row 21
$ws.Range("A21").Value = "103_TruckInsurance_003_InsurantData_002_EnterValuesInWrongFormat"
$ws.Range("B21").Value = "<SET>"
$ws.Range("D21").Value = "Enter values in wrong format"
$ws.Range("H21").Value = "<NOP>"

# Row 22
$ws.Range("A22").Value = "Insurant Data Page check error hint formatting"
$ws.Range("B22").Value = "<CHK>"
$ws.Range("D22").Value = "Check error hint formatting"
$ws.Range("H22").Value = "<NOP>"

# Row 23
$ws.Range("A23").Value = "103_TruckInsurance_003_InsurantData_002_EnterValuesInWrongFormat Part 2"
$ws.Range("B23").Value = "<SET>"
$ws.Range("D23").Value = "Enter values in wrong format part 2"
$ws.Range("H23").Value = "<NOP>"

# Row 24
$ws.Range("A24").Value = "Insurant Data Page check error hint formatting Part 2"
$ws.Range("B24").Value = "<CHK>"
$ws.Range("D24").Value = "Check error hint formatting Part 2"
$ws.Range("H24").Value = "<NOP>"

# Row 25
$ws.Range("A25").Value = "103_TruckInsurance_003_InsurantData_003_ListContents"
$ws.Range("B25").Value = "<CHK>"
$ws.Range("D25").Value = "List content"
$ws.Range("H25").Value = "<NOP>"

# --- Column width adjustments (D:E unified width, G narrowed to match) ---
$ws.Columns.Item(4).ColumnWidth = 38.0
$ws.Columns.Item(5).ColumnWidth = 38.0
$ws.Columns.Item(7).ColumnWidth = 38.0

# --- Reposition the screenshot picture so it sits below the new rows ---
$shp = $ws.Shapes.Item(1)
$shp.Top = 434.4
$shp.Left = 8.4
$shp.Width = 1146.6066929133858
$shp.Height = 719.91

# --- Selection moves to the newly-added last row ---
[void]$ws.Range("D25").Select()
